$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "65.621.38"
$ws.Range("E2").Value = "  -0.06%  "
# Row 3
$ws.Range("D3").Value = "2.655.52"
$ws.Range("E3").Value = "  -0.67%  "
# Row 4
$ws.Range("E4").Value = "  +0.00%  "
# Row 5
$ws.Range("D5").Value = "'596.87"
$ws.Range("E5").Value = "  -0.51%  "
# Row 6
$ws.Range("D6").Value = "'158.54"
$ws.Range("E6").Value = "  +1.37%  "
# Row 7
$ws.Range("B7").Value = "XRP"
$ws.Range("C7").Value = "https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp"
$ws.Range("D7").Value = "'0.641"
$ws.Range("E7").Value = "  +4.42%  "
# Row 8
$ws.Range("B8").Value = "USDC"
$ws.Range("C8").Value = "https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc"
$ws.Range("D8").Value = "'1.00"
$ws.Range("E8").Value = "  -0.02%  "
# Row 9
$ws.Range("E9").Value = "  -1.96%  "
# Row 10
$ws.Range("D10").Value = "'0.398"
$ws.Range("E10").Value = "  -0.21%  "
# Row 11
$ws.Range("D11").Value = "'5.84"
$ws.Range("E11").Value = "  -0.60%  "
# Row 12
$ws.Range("E12").Value = "  +1.07%  "
# Row 13
$ws.Range("D13").Value = "'29.06"
$ws.Range("E13").Value = "  -0.44%  "
# Row 14
$ws.Range("E14").Value = "  -0.98%  "
# Row 15
$ws.Range("D15").Value = "3.131.11"
$ws.Range("E15").Value = "  -0.75%  "
# Row 16
$ws.Range("D16").Value = "65.529.50"
$ws.Range("E16").Value = "  -0.01%  "
# Row 17
$ws.Range("D17").Value = "2.638.44"
$ws.Range("E17").Value = "  -1.42%  "
# Row 18
$ws.Range("D18").Value = "'12.49"
$ws.Range("E18").Value = "  -2.84%  "
# Row 19
$ws.Range("E19").Value = "  -0.59%  "
# Row 20
$ws.Range("D20").Value = "'353.05"
$ws.Range("E20").Value = "  +0.19%  "
# Row 21
$ws.Range("D21").Value = "'7.43"
$ws.Range("E21").Value = "  -1.85%  "
# Row 22
$ws.Range("D22").Value = "'0.999"
$ws.Range("E22").Value = "  -0.05%  "
# Row 23
$ws.Range("D23").Value = "'69.58"
$ws.Range("E23").Value = "  -0.23%  "
# Row 24
$ws.Range("B24").Value = "PEPE"
$ws.Range("C24").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D24").Value = "'0.0000113"
$ws.Range("E24").Value = "  +1.35%  "
# Row 25
$ws.Range("B25").Value = "SuiNetwork"
$ws.Range("C25").Value = "https://coinranking.com/coin/3xJluUMvp+suinetwork-sui"
$ws.Range("D25").Value = "'1.76"
$ws.Range("E25").Value = "  +5.52%  "
# Row 26
$ws.Range("E26").Value = "  +0.52%  "
# Row 27
$ws.Range("E27").Value = "  +1.71%  "
# Row 28
$ws.Range("D28").Value = "'558.95"
$ws.Range("E28").Value = "  +5.89%  "
# Row 29
$ws.Range("B29").Value = "Aptos"
$ws.Range("C29").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D29").Value = "'8.08"
$ws.Range("E29").Value = "  +0.73%  "
# Row 30
$ws.Range("B30").Value = "Kaspa"
$ws.Range("C30").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D30").Value = "'0.163"
$ws.Range("E30").Value = "  -2.23%  "
# Row 31
$ws.Range("E31").Value = "  +0.35%  "
# Row 32
$ws.Range("E32").Value = "  -0.53%  "
# Row 33
$ws.Range("E33").Value = "  +2.14%  "
# Row 34
$ws.Range("D34").Value = "'6.66"
$ws.Range("E34").Value = "  +2.64%  "
# Row 35
$ws.Range("E35").Value = "  -1.14%  "
# Row 36
$ws.Range("E36").Value = "  -0.36%  "
# Row 37
$ws.Range("E37").Value = "  -0.38%  "
# Row 38
$ws.Range("D38").Value = "'0.999"
$ws.Range("E38").Value = "  -0.08%  "
# Row 39
$ws.Range("E39").Value = "  +0.85%  "
# Row 40
$ws.Range("D40").Value = "'153.59"
$ws.Range("E40").Value = "  -2.69%  "
# Row 41
$ws.Range("E41").Value = "  -0.01%  "
# Row 42
$ws.Range("D42").Value = "'2.46"
$ws.Range("E42").Value = "  +5.67%  "
# Row 43
$ws.Range("D43").Value = "'161.05"
$ws.Range("E43").Value = "  -1.41%  "
# Row 44
$ws.Range("E44").Value = "  -1.04%  "
# Row 45
$ws.Range("E45").Value = "  +1.13%  "
# Row 46
$ws.Range("E46").Value = "  +2.24%  "
# Row 47
$ws.Range("D47").Value = "'0.642"
$ws.Range("E47").Value = "  +0.63%  "
# Row 48
$ws.Range("D48").Value = "'0.0257"
$ws.Range("E48").Value = "  +0.05%  "
# Row 49
$ws.Range("E49").Value = "  +1.91%  "
# Row 50
$ws.Range("D50").Value = "'19.73"
$ws.Range("E50").Value = "  -1.99%  "
# Row 51
$ws.Range("D51").Value = "0.0₆0243"
$ws.Range("E51").Value = "  -7.71%  "
